# data base migration from mysql to postgres.
# defect fix: If consumption for current month does not exist, create it.
# -> "Step" sheet gains an Edit-Category automation flow (new OCCURENCE /
#    ELEMENT NO columns + 6 new steps), the "Test Case" sheet selection
#    moves, and the "TC Desription" sheet rows get an explicit height.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Test Case" sheet: just a selection change (G1:G1048576 -> B5)
# ---------------------------------------------------------------------
$wsTestCase = $wb.Worksheets.Item("Test Case")
$wsTestCase.Range("B5").Select()

# ---------------------------------------------------------------------
# 2. "TC Desription" sheet: rows 1 & 2 get an explicit row height
# ---------------------------------------------------------------------
$wsTcDesc = $wb.Worksheets.Item("TC Desription")
$wsTcDesc.Rows.Item(1).RowHeight = 16.5
$wsTcDesc.Rows.Item(2).RowHeight = 16.5

# ---------------------------------------------------------------------
# 3. "Step" sheet: insert two new columns (OCCURENCE / ELEMENT NO)
#    before the old XPATH/VALUE columns, then append the new
#    "Edit Category" test steps (rows 14-19).
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Step")
$ws.Activate()

# Insert 2 columns at F:G -- old F (XPATH) / G (VALUE) shift to H / I
$ws.Range("F1:G1").EntireColumn.Insert()

# New column widths for the freshly inserted F (OCCURENCE) and G (ELEMENT NO)
$ws.Columns.Item(6).ColumnWidth = 9.00500011
$ws.Columns.Item(7).ColumnWidth = 12.63000011

# New column headers
$ws.Range("F1").Value = "OCCURENCE"
$ws.Range("G1").Value = "ELEMENT NO"

# Every existing data row gets "SINGLE" in the new OCCURENCE column
$ws.Range("F2").Value = "SINGLE"
$ws.Range("F3").Value = "SINGLE"
$ws.Range("F4").Value = "SINGLE"
$ws.Range("F5").Value = "SINGLE"
$ws.Range("F6").Value = "SINGLE"
$ws.Range("F7").Value = "SINGLE"
$ws.Range("F8").Value = "SINGLE"
$ws.Range("F9").Value = "SINGLE"
$ws.Range("F10").Value = "SINGLE"
$ws.Range("F11").Value = "SINGLE"
$ws.Range("F12").Value = "SINGLE"
$ws.Range("F13").Value = "SINGLE"

# defect fix: row 10 VALUE (now column I) becomes "Test Category 3"
$ws.Range("I10").Value = "Test Category 3"

# ---- new rows: Edit Category automation ----

# Row 14 - wait for login page to load (re-used step, now also featured
# before the edit-category flow)
$ws.Range("A14").Value = 3
$ws.Range("B14").Value = "Wait for login page to load"
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = "WAIT"
$ws.Range("E14").Value = "PRESENCE"
$ws.Range("F14").Value = "SINGLE"
$ws.Range("H14").Value = "//*[@id='btn-edit-category']"
$ws.Range("I14").Value = "e-Kanban"
$ws.Rows.Item(14).RowHeight = 16.5

# Row 15 - click nth edit button (MULTIPLE occurence, element no 2)
$ws.Range("A15").Value = 3
$ws.Range("B15").Value = "click nth edit button"
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = "CLICK"
$ws.Range("F15").Value = "MULTIPLE"
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = "//*[@id='btn-edit-category']"

# Row 16 - enter new category value
$ws.Range("A16").Value = 3
$ws.Range("B16").Value = "enter new category value"
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = "SENDKEY"
$ws.Range("F16").Value = "SINGLE"
$ws.Range("H16").Value = "//*[@id='form-edit-category']/div[2]/div/span/input[@name='name']"
$ws.Range("I16").Value = "TEST Category 10"

# Row 17 - submit form
$ws.Range("A17").Value = 3
$ws.Range("B17").Value = "Submit form"
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = "SUBMIT"
$ws.Range("F17").Value = "SINGLE"
$ws.Range("H17").Value = "//*[@id='btn-edit-submit']"

# Row 18 - wait for request to complete
$ws.Range("A18").Value = 2
$ws.Range("B18").Value = "wait for request to complete"
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = "WAIT"
$ws.Range("E18").Value = "ABSENCE"
$ws.Range("F18").Value = "SINGLE"
$ws.Range("H18").Value = "//*[@id='form-edit-category']"
$ws.Rows.Item(18).RowHeight = 16.5

# Row 19 - search for updated value
$ws.Range("A19").Value = 3
$ws.Range("B19").Value = "Search for updated value"
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = "SEARCH"
$ws.Range("F19").Value = "SINGLE"
$ws.Range("H19").Value = "//*[@id='list-category']/li"
$ws.Range("I19").Value = "TEST Category 10"

# ---------------------------------------------------------------------
# 4. Rebuild the hyperlinks: the column insert doesn't auto-shift the
#    existing hyperlink anchors, so clear them and re-add on the now
#    correct H-column cells, plus the new edit-category steps.
# ---------------------------------------------------------------------
$ws.Cells.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("I3"), "mailto:zahid7292@gmail.com")
$ws.Hyperlinks.Add($ws.Range("I4"), "mailto:zahid7292@gmail.com")
$ws.Hyperlinks.Add($ws.Range("H2"), "mailto://*[@id='app-heading']")
$ws.Hyperlinks.Add($ws.Range("H6"), "mailto://*[@id='app-header']/div/span")
$ws.Hyperlinks.Add($ws.Range("H7"), "mailto://*[@id='app-header']/nav/a[@href='/profile']")
$ws.Hyperlinks.Add($ws.Range("H8"), "mailto://header[@id='page-header']/div[1]/span")
$ws.Hyperlinks.Add($ws.Range("H9"), "mailto://*[@id='btn-add-category']")
$ws.Hyperlinks.Add($ws.Range("H10"), "mailto://*[@id='form-add-category']/div[2]/div/span/input[@name='name']")
$ws.Hyperlinks.Add($ws.Range("H12"), "mailto://*[@id='form-add-category']")
$ws.Hyperlinks.Add($ws.Range("H5"), "mailto://*[@id='btn-login']")
$ws.Hyperlinks.Add($ws.Range("H13"), "mailto://*[@id='list-category']/li")
$ws.Hyperlinks.Add($ws.Range("H16"), "mailto://*[@id='form-edit-category']/div[2]/div/span/input[@name='name']")
$ws.Hyperlinks.Add($ws.Range("H17"), "mailto://*[@id='btn-edit-submit']")
$ws.Hyperlinks.Add($ws.Range("H19"), "mailto://*[@id='list-category']/li")
$ws.Hyperlinks.Add($ws.Range("H18"), "mailto://*[@id='form-edit-category']")
$ws.Hyperlinks.Add($ws.Range("H14"), "mailto://*[@id='btn-edit-category']")
$ws.Hyperlinks.Add($ws.Range("H15"), "mailto://*[@id='btn-edit-category']")

# ---------------------------------------------------------------------
# 5. Restore the "Step" sheet selection/active state
# ---------------------------------------------------------------------
$ws.Range("E14").Select()
